$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 523.3333
$ws.Range("I9").Value = 523.3333
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 523.3333
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = -354.3333
$ws.Range("N9").Value = $null

$ws.Range("H94").Value = 3053.1428
$ws.Range("I94").Value = 228.66667
$ws.Range("J94").Value = 20000
$ws.Range("K94").Value = 228.66667
$ws.Range("L94").Value = 20000
$ws.Range("M94").Value = 222.33333
$ws.Range("N94").Value = -20902

$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").Value = $null

$ws.Range("H137").Value = 3877.7896
$ws.Range("I137").Value = 4968.4287
$ws.Range("J137").Value = 3241.5833
$ws.Range("K137").Value = 14905.2861
$ws.Range("L137").Value = 9724.749899999999
$ws.Range("M137").Value = -12355.2861
$ws.Range("N137").Value = -14824.7499

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5242.341
$ws.Range("I61").Value = 2178.5625
$ws.Range("K61").Value = 2178.5625
$ws.Range("M61").Value = -1966.5625

$ws.Range("H74").Value = 46942.406
$ws.Range("I74").Value = 82255.75
$ws.Range("J74").Value = 5397.294
$ws.Range("K74").Value = 82255.75
$ws.Range("L74").Value = 5397.294
$ws.Range("M74").Value = -81381.75
$ws.Range("N74").Value = -7145.294

$ws.Range("H77").Value = 46942.406
$ws.Range("I77").Value = 82255.75
$ws.Range("J77").Value = 5397.294
$ws.Range("K77").Value = 411278.75
$ws.Range("L77").Value = 26986.47
$ws.Range("M77").Value = -406910.75
$ws.Range("N77").Value = -35722.47

$ws.Range("H97").Value = 16698989
$ws.Range("I97").Value = 799.5
$ws.Range("K97").Value = 799.5
$ws.Range("M97").Value = -303.5

$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").Value = $null

$ws.Range("H122").Value = 15937.375
$ws.Range("I122").Value = 21449.8
$ws.Range("K122").Value = 64349.39999999999
$ws.Range("M122").Value = -61899.39999999999

$ws.Range("H132").Value = 6887.838
$ws.Range("I132").Value = 5504
$ws.Range("J132").Value = 9770.833000000001
$ws.Range("K132").Value = 16512
$ws.Range("L132").Value = 29312.499
$ws.Range("M132").Value = -13982
$ws.Range("N132").Value = -34372.499

$ws.Range("H136").Value = 5242.341
$ws.Range("I136").Value = 2178.5625
$ws.Range("K136").Value = 6535.6875
$ws.Range("M136").Value = -3985.6875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 27778060
$ws.Range("I80").Value = 62500296
$ws.Range("J80").Value = 273
$ws.Range("K80").Value = 62500296
$ws.Range("L80").Value = 273
$ws.Range("M80").Value = -62499298
$ws.Range("N80").Value = -2269

$ws.Range("H83").Value = 27778060
$ws.Range("I83").Value = 62500296
$ws.Range("J83").Value = 273
$ws.Range("K83").Value = 312501480
$ws.Range("L83").Value = 1365
$ws.Range("M83").Value = -312496488
$ws.Range("N83").Value = -11349

$ws.Range("H105").Value = 88296.89
$ws.Range("I105").Value = 104990
$ws.Range("K105").Value = 104990
$ws.Range("M105").Value = -103243

$ws.Range("H108").Value = 59376
$ws.Range("J108").Value = 59376
$ws.Range("L108").Value = 59376
$ws.Range("N108").Value = -67056

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 13339554
$ws.Range("I132").Value = 2670.111
$ws.Range("K132").Value = 8010.333
$ws.Range("M132").Value = -5480.333

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 2759
$ws.Range("I3").Value = 2759
$ws.Range("K3").Value = 8277
$ws.Range("M3").Value = -8165

$ws.Range("H4").Value = 224299090
$ws.Range("J4").Value = 168114980
$ws.Range("L4").Value = 504344940
$ws.Range("N4").Value = -504345164

$ws.Range("H38").Value = 65.333336
$ws.Range("I38").Value = 64
$ws.Range("K38").Value = 192
$ws.Range("M38").Value = 155

$ws.Range("H56").Value = 7057.5454
$ws.Range("I56").Value = 7057.5454
$ws.Range("K56").Value = 7057.5454
$ws.Range("M56").Value = -6527.5454

$ws.Range("H114").Value = 544.1177
$ws.Range("I114").Value = 380.83334
$ws.Range("J114").Value = 633.1818
$ws.Range("K114").Value = 1142.50002
$ws.Range("L114").Value = 1899.5454
$ws.Range("M114").Value = 2111.49998
$ws.Range("N114").Value = -8407.545399999999

$ws.Range("H120").Value = 0
$ws.Range("I120").Value = 0
$ws.Range("K120").Value = 0
$ws.Range("M120").Value = $null

$ws.Range("H122").Value = 4044475.2
$ws.Range("I122").Value = 5659265.5
$ws.Range("K122").Value = 50933389.5
$ws.Range("M122").Value = -50930939.5

$ws.Range("H132").Value = 13434.823
$ws.Range("I132").Value = 5833.3335
$ws.Range("J132").Value = 21986.5
$ws.Range("K132").Value = 52500.0015
$ws.Range("L132").Value = 197878.5
$ws.Range("M132").Value = -49970.0015
$ws.Range("N132").Value = -202938.5

$ws.Range("H136").Value = 1417.5
$ws.Range("I136").Value = 905.7143
$ws.Range("K136").Value = 2717.1429
$ws.Range("M136").Value = 2382.8571

$ws.Range("H140").Value = 107500.9
$ws.Range("I140").Value = 134601.4
$ws.Range("J140").Value = 5874
$ws.Range("K140").Value = 403804.2
$ws.Range("L140").Value = 17622
$ws.Range("M140").Value = -398624.2
$ws.Range("N140").Value = -27982

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5413.2856
$ws.Range("I80").Value = 3999.5
$ws.Range("J80").Value = 5978.8
$ws.Range("K80").Value = 3999.5
$ws.Range("L80").Value = 5978.8
$ws.Range("M80").Value = -3001.5
$ws.Range("N80").Value = -7974.8

$ws.Range("H83").Value = 5413.2856
$ws.Range("I83").Value = 3999.5
$ws.Range("J83").Value = 5978.8
$ws.Range("K83").Value = 19997.5
$ws.Range("L83").Value = 29894
$ws.Range("M83").Value = -15005.5
$ws.Range("N83").Value = -39878

$ws.Range("H122").Value = 2266084.2
$ws.Range("I122").Value = 3294147.2
$ws.Range("K122").Value = 9882441.600000001
$ws.Range("M122").Value = -9879991.600000001

$ws.Range("H126").Value = 2804.5833
$ws.Range("I126").Value = 2834.6365
$ws.Range("K126").Value = 8503.9095
$ws.Range("M126").Value = -6033.9095

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5008.4614
$ws.Range("I7").Value = 3400.8333
$ws.Range("K7").Value = 3400.8333
$ws.Range("M7").Value = -3288.8333

$ws.Range("H22").Value = 750.25
$ws.Range("I22").Value = 667
$ws.Range("K22").Value = 667
$ws.Range("M22").Value = -372

$ws.Range("H27").Value = 750.25
$ws.Range("I27").Value = 667
$ws.Range("K27").Value = 667
$ws.Range("M27").Value = -560

$ws.Range("H36").Value = 51598
$ws.Range("J36").Value = 51598
$ws.Range("L36").Value = 51598
$ws.Range("N36").Value = -52722

$ws.Range("H40").Value = 3949.0303
$ws.Range("I40").Value = 2950.9285
$ws.Range("K40").Value = 2950.9285
$ws.Range("M40").Value = -2814.9285

$ws.Range("H68").Value = 3723.0588
$ws.Range("I68").Value = 1928.9
$ws.Range("K68").Value = 1928.9
$ws.Range("M68").Value = -1179.9

$ws.Range("H71").Value = 3723.0588
$ws.Range("I71").Value = 1928.9
$ws.Range("K71").Value = 9644.5
$ws.Range("M71").Value = -5900.5

$ws.Range("H82").Value = 1835.75
$ws.Range("J82").Value = 2672
$ws.Range("L82").Value = 2672
$ws.Range("N82").Value = -3394

$ws.Range("H85").Value = 1835.75
$ws.Range("J85").Value = 2672
$ws.Range("L85").Value = 2672
$ws.Range("N85").Value = -5168

$ws.Range("H93").Value = 12432.556
$ws.Range("I93").Value = 10332.167
$ws.Range("K93").Value = 10332.167
$ws.Range("M93").Value = -9084.166999999999

$ws.Range("H126").Value = 5008.4614
$ws.Range("I126").Value = 3400.8333
$ws.Range("K126").Value = 10202.4999
$ws.Range("M126").Value = -7732.499899999999

$ws.Range("H132").Value = 17866786
$ws.Range("I132").Value = 71432856
$ws.Range("J132").Value = 11428.523
$ws.Range("K132").Value = 214298568
$ws.Range("L132").Value = 34285.569
$ws.Range("M132").Value = -214296038
$ws.Range("N132").Value = -39345.569

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 16097
$ws.Range("I54").Value = 14535.8
$ws.Range("K54").Value = 14535.8
$ws.Range("M54").Value = -14015.8

$ws.Range("H124").Value = 29084.8
$ws.Range("I124").Value = 15000
$ws.Range("J124").Value = 32606
$ws.Range("K124").Value = 15000
$ws.Range("L124").Value = 32606
$ws.Range("M124").Value = -10090
$ws.Range("N124").Value = -42426

$ws.Range("H126").Value = 7577
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 7577
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 22731
$ws.Range("M126").Value = $null
$ws.Range("N126").Value = -27671

$ws.Range("H132").Value = 22745712
$ws.Range("I132").Value = 45464572
$ws.Range("J132").Value = 26851.818
$ws.Range("K132").Value = 136393716
$ws.Range("L132").Value = 80555.454
$ws.Range("M132").Value = -136391186
$ws.Range("N132").Value = -85615.454

$ws.Range("H136").Value = 23283630
$ws.Range("I136").Value = 45455532
$ws.Range("K136").Value = 136366596
$ws.Range("M136").Value = -136364046
